$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $value) {
    $range = $ws.Range($ref)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws "D2" "65.088.39"
$ws.Range("E2").Value = "  +0.58%  "
Set-TextCell $ws "D3" "3.539.47"
$ws.Range("E3").Value = "  +4.59%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextCell $ws "D5" "598.48"
$ws.Range("E5").Value = "  +3.90%  "
Set-TextCell $ws "D6" "138.01"
$ws.Range("E6").Value = "  +3.12%  "
Set-TextCell $ws "D7" "3.538.94"
$ws.Range("E7").Value = "  +4.58%  "
$ws.Range("E8").Value = "  +0.16%  "
Set-TextCell $ws "D9" "0.494"
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("E10").Value = "  +4.15%  "
Set-TextCell $ws "D11" "6.92"
$ws.Range("E11").Value = "  -1.32%  "
$ws.Range("E12").Value = "  +4.81%  "
Set-TextCell $ws "D13" "4.141.91"
$ws.Range("E13").Value = "  +4.64%  "
Set-TextCell $ws "D15" "27.38"
$ws.Range("E15").Value = "  +6.00%  "
Set-TextCell $ws "D16" "3.541.17"
$ws.Range("E16").Value = "  +3.54%  "
$ws.Range("E17").Value = "  +1.74%  "
Set-TextCell $ws "D18" "65.088.52"
$ws.Range("E18").Value = "  +0.59%  "
Set-TextCell $ws "D19" "10.08"
$ws.Range("E19").Value = "  +6.88%  "
$ws.Range("E20").Value = "  +2.41%  "
Set-TextCell $ws "D21" "14.22"
$ws.Range("E21").Value = "  +6.39%  "
Set-TextCell $ws "D22" "391.93"
$ws.Range("E22").Value = "  +3.93%  "
$ws.Range("E23").Value = "  +5.43%  "
Set-TextCell $ws "D24" "3.682.61"
$ws.Range("E24").Value = "  +4.59%  "
Set-TextCell $ws "D25" "73.85"
$ws.Range("E25").Value = "  +3.27%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +11.23%  "
Set-TextCell $ws "D28" "7.82"
$ws.Range("E28").Value = "  +13.24%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell $ws "D29" "2.27"
$ws.Range("E29").Value = "  +5.19%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws "D30" "8.30"
$ws.Range("E30").Value = "  +5.39%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextCell $ws "D31" "0.970"
$ws.Range("E31").Value = "  -3.18%  "
Set-TextCell $ws "D32" "3.558.12"
$ws.Range("E32").Value = "  +4.77%  "
Set-TextCell $ws "D33" "1.40"
$ws.Range("E33").Value = "  +24.89%  "
$ws.Range("E34").Value = "  +0.01%  "
Set-TextCell $ws "D35" "23.86"
$ws.Range("E35").Value = "  +5.37%  "
Set-TextCell $ws "D36" "0.145"
$ws.Range("E36").Value = "  +3.34%  "
$ws.Range("E37").Value = "  +10.98%  "
Set-TextCell $ws "D38" "169.57"
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("E39").Value = "  +5.89%  "
Set-TextCell $ws "D40" "5.03"
$ws.Range("E40").Value = "  +9.68%  "
Set-TextCell $ws "D41" "0.0805"
$ws.Range("E41").Value = "  +8.69%  "
Set-TextCell $ws "D42" "0.823"
$ws.Range("E42").Value = "  +2.66%  "
$ws.Range("E43").Value = "  +23.20%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell $ws "D44" "1.00"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws "D45" "42.40"
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("E46").Value = "  +3.99%  "
$ws.Range("E47").Value = "  +12.39%  "
Set-TextCell $ws "D48" "1.68"
$ws.Range("E48").Value = "  +6.77%  "
$ws.Range("E49").Value = "  +7.35%  "
Set-TextCell $ws "D50" "2.395.14"
$ws.Range("E50").Value = "  +11.91%  "
Set-TextCell $ws "D51" "310.45"
$ws.Range("E51").Value = "  +19.04%  "
